# Auto-generated edit script applying the cryptos.xlsx price/volume/hour refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin name swap + link swap between rows 7 and 8) ---
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"

# --- Numeric-looking cells that must stay stored as text (price, % change, hour) ---
# Force text storage (NumberFormat "@") so Excel does not coerce these into numbers,
# then restore the default "Normal" style so no stray style index is left on the cell.
$forceTextCells = @(
    "D2", "E2", "G2", "D3", "E3", "G3", "D4", "E4", "G4", "D5",
    "E5", "G5", "D6", "E6", "G6", "D7", "E7", "G7", "D8", "E8",
    "G8", "E9", "G9", "D10", "E10", "G10", "D11", "E11", "G11", "D12",
    "E12", "G12", "D13", "E13", "G13", "D14", "E14", "G14", "D15", "E15",
    "G15", "D16", "E16", "G16", "D17", "E17", "G17", "D18", "E18", "G18",
    "D19", "E19", "G19", "D20", "E20", "G20", "E21", "G21", "D22", "E22",
    "G22", "D23", "E23", "G23", "D24", "E24", "G24", "D25", "E25", "G25",
    "D26", "E26", "G26", "D27", "E27", "G27", "G28", "G29", "G30", "G31",
    "G32", "G33", "G34", "G35", "G36", "G37", "G38", "D39", "E39", "G39",
    "D40", "E40", "G40", "D41", "E41", "G41", "D42", "E42", "G42", "D43",
    "E43", "G43", "D44", "E44", "G44", "D45", "E45", "G45", "D46", "E46",
    "G46", "D47", "E47", "G47", "D48", "E48", "G48", "D49", "E49", "G49",
    "D50", "E50", "G50", "D51", "E51", "G51"
)
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "311.10"
$ws.Range("E2").Value = "0.39%"
$ws.Range("G2").Value = "9"
$ws.Range("D3").Value = "37.68"
$ws.Range("E3").Value = "-1.42%"
$ws.Range("G3").Value = "9"
$ws.Range("D4").Value = "5.089"
$ws.Range("E4").Value = "-0.95%"
$ws.Range("G4").Value = "9"
$ws.Range("D5").Value = "0.07771"
$ws.Range("E5").Value = "-2.89%"
$ws.Range("G5").Value = "9"
$ws.Range("D6").Value = "4.367"
$ws.Range("E6").Value = "-2.47%"
$ws.Range("G6").Value = "9"
$ws.Range("D7").Value = "1.896"
$ws.Range("E7").Value = "-8.32%"
$ws.Range("G7").Value = "9"
$ws.Range("D8").Value = "8.217"
$ws.Range("E8").Value = "-1.09%"
$ws.Range("G8").Value = "9"
$ws.Range("E9").Value = "-7.26%"
$ws.Range("G9").Value = "9"
$ws.Range("D10").Value = "0.9213"
$ws.Range("E10").Value = "-2.06%"
$ws.Range("G10").Value = "9"
$ws.Range("D11").Value = "0.1215"
$ws.Range("E11").Value = "-6.85%"
$ws.Range("G11").Value = "9"
$ws.Range("D12").Value = "0.1921"
$ws.Range("E12").Value = "-0.07%"
$ws.Range("G12").Value = "9"
$ws.Range("D13").Value = "0.09337"
$ws.Range("E13").Value = "4.14%"
$ws.Range("G13").Value = "9"
$ws.Range("D14").Value = "0.03420"
$ws.Range("E14").Value = "-1.63%"
$ws.Range("G14").Value = "9"
$ws.Range("D15").Value = "0.09684"
$ws.Range("E15").Value = "-0.36%"
$ws.Range("G15").Value = "9"
$ws.Range("D16").Value = "0.001367"
$ws.Range("E16").Value = "-3.37%"
$ws.Range("G16").Value = "9"
$ws.Range("D17").Value = "0.005960"
$ws.Range("E17").Value = "-11.80%"
$ws.Range("G17").Value = "9"
$ws.Range("D18").Value = "3.550"
$ws.Range("E18").Value = "-0.93%"
$ws.Range("G18").Value = "9"
$ws.Range("D19").Value = "0.3400"
$ws.Range("E19").Value = "-1.88%"
$ws.Range("G19").Value = "9"
$ws.Range("D20").Value = "5.262"
$ws.Range("E20").Value = "4.37%"
$ws.Range("G20").Value = "9"
$ws.Range("E21").Value = "0.16%"
$ws.Range("G21").Value = "9"
$ws.Range("D22").Value = "0.2591"
$ws.Range("E22").Value = "2.33%"
$ws.Range("G22").Value = "9"
$ws.Range("D23").Value = "0.02104"
$ws.Range("E23").Value = "5,587.57%"
$ws.Range("G23").Value = "9"
$ws.Range("D24").Value = "0.04345"
$ws.Range("E24").Value = "-0.70%"
$ws.Range("G24").Value = "9"
$ws.Range("D25").Value = "0.001213"
$ws.Range("E25").Value = "-2.93%"
$ws.Range("G25").Value = "9"
$ws.Range("D26").Value = "0.004257"
$ws.Range("E26").Value = "-9.09%"
$ws.Range("G26").Value = "9"
$ws.Range("D27").Value = "0.0001300"
$ws.Range("E27").Value = "-63.77%"
$ws.Range("G27").Value = "9"
$ws.Range("G28").Value = "9"
$ws.Range("G29").Value = "9"
$ws.Range("G30").Value = "9"
$ws.Range("G31").Value = "9"
$ws.Range("G32").Value = "9"
$ws.Range("G33").Value = "9"
$ws.Range("G34").Value = "9"
$ws.Range("G35").Value = "9"
$ws.Range("G36").Value = "9"
$ws.Range("G37").Value = "9"
$ws.Range("G38").Value = "9"
$ws.Range("D39").Value = "0.02104"
$ws.Range("E39").Value = "-3.49%"
$ws.Range("G39").Value = "9"
$ws.Range("D40").Value = "0.05024"
$ws.Range("E40").Value = "-2.78%"
$ws.Range("G40").Value = "9"
$ws.Range("D41").Value = "0.007671"
$ws.Range("E41").Value = "0.43%"
$ws.Range("G41").Value = "9"
$ws.Range("D42").Value = "0.009827"
$ws.Range("E42").Value = "-1.71%"
$ws.Range("G42").Value = "9"
$ws.Range("D43").Value = "0.1343"
$ws.Range("E43").Value = "-3.29%"
$ws.Range("G43").Value = "9"
$ws.Range("D44").Value = "0.002060"
$ws.Range("E44").Value = "0.72%"
$ws.Range("G44").Value = "9"
$ws.Range("D45").Value = "0.008874"
$ws.Range("E45").Value = "-2.77%"
$ws.Range("G45").Value = "9"
$ws.Range("D46").Value = "0.00006661"
$ws.Range("E46").Value = "-0.32%"
$ws.Range("G46").Value = "9"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.79%"
$ws.Range("G47").Value = "9"
$ws.Range("D48").Value = "0.002935"
$ws.Range("E48").Value = "-2.96%"
$ws.Range("G48").Value = "9"
$ws.Range("D49").Value = "0.001200"
$ws.Range("E49").Value = "-0.84%"
$ws.Range("G49").Value = "9"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "-0.79%"
$ws.Range("G50").Value = "9"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "-0.79%"
$ws.Range("G51").Value = "9"

foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = "Normal"
}
